# Update the two-digit ÷ one-digit division problems in the table.
# Each data row (1, 5, 9, 13, 17) of the single table holds 5 problems,
# with blank rows in between for pupils to write the answer.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="66÷7="},
    @{Row=1;  Col=2; New="64÷6="},
    @{Row=1;  Col=3; New="99÷8="},
    @{Row=1;  Col=4; New="19÷4="},
    @{Row=1;  Col=5; New="58÷3="},

    @{Row=5;  Col=1; New="64÷3="},
    @{Row=5;  Col=2; New="29÷7="},
    @{Row=5;  Col=3; New="45÷8="},
    @{Row=5;  Col=4; New="40÷3="},
    @{Row=5;  Col=5; New="17÷4="},

    @{Row=9;  Col=1; New="96÷9="},
    @{Row=9;  Col=2; New="54÷2="},
    @{Row=9;  Col=3; New="15÷5="},
    @{Row=9;  Col=4; New="54÷9="},
    @{Row=9;  Col=5; New="82÷3="},

    @{Row=13; Col=1; New="17÷6="},
    @{Row=13; Col=2; New="54÷3="},
    @{Row=13; Col=3; New="66÷2="},
    @{Row=13; Col=4; New="70÷2="},
    @{Row=13; Col=5; New="25÷9="},

    @{Row=17; Col=1; New="68÷4="},
    @{Row=17; Col=2; New="18÷6="},
    @{Row=17; Col=3; New="94÷3="},
    @{Row=17; Col=4; New="53÷3="},
    @{Row=17; Col=5; New="17÷5="}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark/paragraph-mark characters Word appends
    # to cell ranges, keeping only the visible text portion.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $u.New
}
